$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Update Trials sheet rows 302-401 (trials 300-399) with new Iterations/Runtime/Success values
$ws1.Range("B302").Value = 7875
$ws1.Range("C302").Value = 7.914174795150757
$ws1.Range("B303").Value = 1992
$ws1.Range("C303").Value = 1.177507400512695
$ws1.Range("B304").Value = 2325
$ws1.Range("C304").Value = 1.252694129943848
$ws1.Range("B305").Value = 1545
$ws1.Range("C305").Value = 0.6562430858612061
$ws1.Range("B306").Value = 6373
$ws1.Range("C306").Value = 6.542188882827759
$ws1.Range("B307").Value = 1395
$ws1.Range("C307").Value = 0.5318694114685059
$ws1.Range("B308").Value = 2483
$ws1.Range("C308").Value = 1.53418493270874
$ws1.Range("B309").Value = 8012
$ws1.Range("C309").Value = 7.217840433120728
$ws1.Range("B310").Value = 991
$ws1.Range("C310").Value = 0.3791201114654541
$ws1.Range("D310").Value = 1
$ws1.Range("B311").Value = 3278
$ws1.Range("C311").Value = 2.209597826004028
$ws1.Range("B312").Value = 5043
$ws1.Range("C312").Value = 3.890109062194824
$ws1.Range("B313").Value = 15349
$ws1.Range("C313").Value = 21.41181588172913
$ws1.Range("B314").Value = 4893
$ws1.Range("C314").Value = 3.987114191055298
$ws1.Range("B315").Value = 10091
$ws1.Range("C315").Value = 10.00911092758179
$ws1.Range("B316").Value = 11580
$ws1.Range("C316").Value = 13.57191133499146
$ws1.Range("B317").Value = 8226
$ws1.Range("C317").Value = 8.061393976211548
$ws1.Range("B318").Value = 13535
$ws1.Range("C318").Value = 20.2875030040741
$ws1.Range("B319").Value = 7497
$ws1.Range("C319").Value = 6.624865770339966
$ws1.Range("B320").Value = 8761
$ws1.Range("C320").Value = 8.537094831466675
$ws1.Range("B321").Value = 14242
$ws1.Range("C321").Value = 21.55966925621033
$ws1.Range("B322").Value = 2234
$ws1.Range("C322").Value = 1.146967649459839
$ws1.Range("B323").Value = 5640
$ws1.Range("C323").Value = 3.964053153991699
$ws1.Range("B324").Value = 2687
$ws1.Range("C324").Value = 1.48234224319458
$ws1.Range("B325").Value = 8255
$ws1.Range("C325").Value = 7.792657852172852
$ws1.Range("D325").Value = 1
$ws1.Range("B326").Value = 7685
$ws1.Range("C326").Value = 7.290755748748779
$ws1.Range("B327").Value = 2654
$ws1.Range("C327").Value = 1.493052959442139
$ws1.Range("B328").Value = 1177
$ws1.Range("C328").Value = 0.5514485836029053
$ws1.Range("B329").Value = 10092
$ws1.Range("C329").Value = 10.46492910385132
$ws1.Range("B330").Value = 109
$ws1.Range("C330").Value = 0.04421353340148926
$ws1.Range("B331").Value = 3363
$ws1.Range("C331").Value = 1.750082731246948
$ws1.Range("B332").Value = 26608
$ws1.Range("C332").Value = 60.00157880783081
$ws1.Range("D332").Value = 0
$ws1.Range("B333").Value = 4802
$ws1.Range("C333").Value = 3.112475156784058
$ws1.Range("B334").Value = 3908
$ws1.Range("C334").Value = 2.493205785751343
$ws1.Range("B335").Value = 3327
$ws1.Range("C335").Value = 1.775290966033936
$ws1.Range("B336").Value = 1411
$ws1.Range("C336").Value = 0.585538387298584
$ws1.Range("B337").Value = 767
$ws1.Range("C337").Value = 0.2513682842254639
$ws1.Range("B338").Value = 7409
$ws1.Range("C338").Value = 5.296584844589233
$ws1.Range("B339").Value = 3525
$ws1.Range("C339").Value = 0.8395931720733643
$ws1.Range("B340").Value = 5568
$ws1.Range("C340").Value = 3.829446315765381
$ws1.Range("B341").Value = 5988
$ws1.Range("C341").Value = 4.176910400390625
$ws1.Range("B342").Value = 19066
$ws1.Range("C342").Value = 35.93848395347595
$ws1.Range("B343").Value = 4771
$ws1.Range("C343").Value = 3.122979640960693
$ws1.Range("B344").Value = 10717
$ws1.Range("C344").Value = 12.06258964538574
$ws1.Range("B345").Value = 1912
$ws1.Range("C345").Value = 0.8826680183410645
$ws1.Range("B346").Value = 8821
$ws1.Range("C346").Value = 8.9626784324646
$ws1.Range("D346").Value = 1
$ws1.Range("B347").Value = 8140
$ws1.Range("C347").Value = 5.802713871002197
$ws1.Range("B348").Value = 8334
$ws1.Range("C348").Value = 7.182803630828857
$ws1.Range("B349").Value = 5072
$ws1.Range("C349").Value = 2.673810958862305
$ws1.Range("B350").Value = 1622
$ws1.Range("C350").Value = 0.7953050136566162
$ws1.Range("B351").Value = 3809
$ws1.Range("C351").Value = 1.570146799087524
$ws1.Range("B352").Value = 2818
$ws1.Range("C352").Value = 1.459342956542969
$ws1.Range("B353").Value = 3726
$ws1.Range("C353").Value = 2.080800533294678
$ws1.Range("B354").Value = 11179
$ws1.Range("C354").Value = 14.1690034866333
$ws1.Range("B355").Value = 2381
$ws1.Range("C355").Value = 0.9512910842895508
$ws1.Range("B356").Value = 9612
$ws1.Range("C356").Value = 10.42948961257935
$ws1.Range("D356").Value = 1
$ws1.Range("B357").Value = 534
$ws1.Range("C357").Value = 0.1820292472839355
$ws1.Range("B358").Value = 11967
$ws1.Range("C358").Value = 15.67616534233093
$ws1.Range("B359").Value = 24782
$ws1.Range("C359").Value = 60.00254464149475
$ws1.Range("B360").Value = 3062
$ws1.Range("C360").Value = 1.770901203155518
$ws1.Range("B361").Value = 695
$ws1.Range("C361").Value = 0.27638840675354
$ws1.Range("B362").Value = 6960
$ws1.Range("C362").Value = 6.251266717910767
$ws1.Range("B363").Value = 585
$ws1.Range("C363").Value = 0.1995251178741455
$ws1.Range("B364").Value = 10370
$ws1.Range("C364").Value = 11.73479962348938
$ws1.Range("B365").Value = 2220
$ws1.Range("C365").Value = 1.128860235214233
$ws1.Range("B366").Value = 623
$ws1.Range("C366").Value = 0.2330079078674316
$ws1.Range("B367").Value = 1098
$ws1.Range("C367").Value = 0.46468186378479
$ws1.Range("B368").Value = 794
$ws1.Range("C368").Value = 0.2882580757141113
$ws1.Range("B369").Value = 4404
$ws1.Range("C369").Value = 2.968374729156494
$ws1.Range("B370").Value = 1041
$ws1.Range("C370").Value = 0.4369630813598633
$ws1.Range("B371").Value = 796
$ws1.Range("C371").Value = 0.3155262470245361
$ws1.Range("B372").Value = 15674
$ws1.Range("C372").Value = 23.11936092376709
$ws1.Range("B373").Value = 3125
$ws1.Range("C373").Value = 1.590145111083984
$ws1.Range("D373").Value = 1
$ws1.Range("B374").Value = 5201
$ws1.Range("C374").Value = 3.744139909744263
$ws1.Range("B375").Value = 14225
$ws1.Range("C375").Value = 20.66100978851318
$ws1.Range("D375").Value = 1
$ws1.Range("B376").Value = 199
$ws1.Range("C376").Value = 0.07381200790405273
$ws1.Range("B377").Value = 1514
$ws1.Range("C377").Value = 0.7251045703887939
$ws1.Range("B378").Value = 14134
$ws1.Range("C378").Value = 20.41955208778381
$ws1.Range("B379").Value = 3161
$ws1.Range("C379").Value = 2.400487422943115
$ws1.Range("B380").Value = 8951
$ws1.Range("C380").Value = 8.994839191436768
$ws1.Range("B381").Value = 5113
$ws1.Range("C381").Value = 2.06085991859436
$ws1.Range("B382").Value = 19374
$ws1.Range("C382").Value = 33.59503054618835
$ws1.Range("B383").Value = 20212
$ws1.Range("C383").Value = 40.84259605407715
$ws1.Range("B384").Value = 9252
$ws1.Range("C384").Value = 8.95359468460083
$ws1.Range("B385").Value = 1906
$ws1.Range("C385").Value = 0.8900189399719238
$ws1.Range("D385").Value = 1
$ws1.Range("B386").Value = 4621
$ws1.Range("C386").Value = 2.286929845809937
$ws1.Range("B387").Value = 4470
$ws1.Range("C387").Value = 2.687897682189941
$ws1.Range("D387").Value = 1
$ws1.Range("B388").Value = 9649
$ws1.Range("C388").Value = 9.763243675231934
$ws1.Range("B389").Value = 25387
$ws1.Range("C389").Value = 60.00049448013306
$ws1.Range("D389").Value = 0
$ws1.Range("B390").Value = 12007
$ws1.Range("C390").Value = 13.79620456695557
$ws1.Range("D390").Value = 1
$ws1.Range("B391").Value = 12388
$ws1.Range("C391").Value = 13.89958143234253
$ws1.Range("B392").Value = 8076
$ws1.Range("C392").Value = 6.256377935409546
$ws1.Range("B393").Value = 13788
$ws1.Range("C393").Value = 15.16164970397949
$ws1.Range("B394").Value = 3642
$ws1.Range("C394").Value = 1.658358812332153
$ws1.Range("D394").Value = 1
$ws1.Range("B395").Value = 7256
$ws1.Range("C395").Value = 4.850375652313232
$ws1.Range("B396").Value = 682
$ws1.Range("C396").Value = 0.2005040645599365
$ws1.Range("B397").Value = 6265
$ws1.Range("C397").Value = 3.536361217498779
$ws1.Range("D397").Value = 1
$ws1.Range("B398").Value = 3626
$ws1.Range("C398").Value = 1.842355966567993
$ws1.Range("B399").Value = 1688
$ws1.Range("C399").Value = 0.6145164966583252
$ws1.Range("B400").Value = 2388
$ws1.Range("C400").Value = 0.9197437763214111
$ws1.Range("B401").Value = 1836
$ws1.Range("C401").Value = 0.562293529510498
$ws1.Range("D401").Value = 1
# Update Summary sheet Success Ratio
$ws2.Range("C2").Value = 0.2425
